$d = $word.ActiveDocument

# --- Edit 1: "judge системата" -> split into "judge" + " системата" with proofErr spellStart/spellEnd ---
$p1 = $d.Paragraphs(3).Range
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0BC86920" w14:textId="77777777"><w:pPr><w:spacing w:before="40" w:after="40"/><w:jc w:val="center"/><w:rPr><w:lang w:val="bg-BG"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>Тествайте</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> решенията си в </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>judge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> системата</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r></w:p>
'@
$p1.InsertXML($xml1)

# --- Edit 2: pool sentence - split "запълненост ..." / "... тръба" runs, add spell/gram proofErr marks ---
$p2 = $d.Paragraphs(16).Range
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="53CF6157" w14:textId="5E55E74A"><w:pPr><w:pStyle w:val="ac"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="32"/></w:numPr><w:rPr><w:lang w:val="bg-BG"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>"</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/></w:rPr><w:t xml:space="preserve">The pool is </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>запълненост</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> на басейна в проценти</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/></w:rPr><w:t xml:space="preserve">% full. Pipe 1: </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">процент вода от първата </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>тръба</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/></w:rPr><w:t>%</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/></w:rPr><w:t xml:space="preserve">. Pipe 2: </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">процент вода от втората </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>тръба</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/></w:rPr><w:t>%</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:b/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>"</w:t></w:r></w:p>
'@
$p2.InsertXML($xml2)

# --- Edit 3: fuel paragraph - split off a run of 5 spaces before "В случай..." ---
$p3 = $d.Paragraphs(266).Range
$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="35FB6CCF" w14:textId="7928840A"><w:pPr><w:rPr><w:lang w:val="bg-BG"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">Напишете програма, която познава дали резервоара на едно превозно средство има нужда от презареждане на горивото или не. От конзолата се четат </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>два реда</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>текст</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> и реално число</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>, на първия ред се чете типа на горивото</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> – текст с възможности:</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:lang w:val="bg-BG"/></w:rPr><w:t>"</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:b/></w:rPr><w:t>Diesel</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:lang w:val="bg-BG"/></w:rPr><w:t>"</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:lang w:val="bg-BG"/></w:rPr><w:t>"</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:b/></w:rPr><w:t>Gasoline</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:lang w:val="bg-BG"/></w:rPr><w:t>"</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="bg-BG"/></w:rPr><w:t>или</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:lang w:val="bg-BG"/></w:rPr><w:t>"</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:b/></w:rPr><w:t>Gas</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:lang w:val="bg-BG"/></w:rPr><w:t>"</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>а на втория литрите гориво</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> които има в резервоара</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">. Ако литрите гориво </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>са повече или равни на 25</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>, на конзолата да се отпечата "</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:b/></w:rPr><w:t xml:space="preserve">You have enough </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>вида на горивото</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:b/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>"</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">ако </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>са по-малко от 25</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>, да се отпечата "</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:b/></w:rPr><w:t xml:space="preserve">Fill your tank with </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>вида на горивото</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:b/></w:rPr><w:t>!</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">". </w:t></w:r><w:r><w:t xml:space="preserve">     </w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">В случай, че бъде въведено гориво, </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>различно от посоченото</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>да се отпечата "</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:b/></w:rPr><w:t>Invalid fuel!</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>".</w:t></w:r></w:p>
'@
$p3.InsertXML($xml3)
